$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "NA" value from C49 (it moves down to the new row 50). Setting a
# lone apostrophe stores a genuine empty TEXT value (matching how the other
# "no page number" rows in column C are stored) instead of blanking the cell
# out entirely; resetting the style afterwards drops the quote-prefix
# formatting flag so no stray formatting is left on the cell.
$ws.Range("C49").Value = "'"
$ws.Range("C49").Style = "Normal"

# Add the new row 50 with the data that was previously on row 49's C column.
# Prefix the date with an apostrophe so Excel stores it as literal text
# rather than auto-converting it to a numeric date serial value, then reset
# the cell style back to Normal so no stray formatting is left behind.
$ws.Range("A50").Value = "'2025-04-15"
$ws.Range("A50").Style = "Normal"
$ws.Range("B50").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C50").Value = "NA"
$ws.Range("D50").Value = 1
